$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Row 1: "100" -> "0M"
$t.Cell(1,1).Range.Text = "0M"

# 2) Row 2: "0" -> "0M"
$t.Cell(2,1).Range.Text = "0M"

# 3) Row 3: "63" -> "0M"
$t.Cell(3,1).Range.Text = "0M"

# 4) Insert 10 new rows right after row 3 (before the old row 4),
#    each containing a single value.
$newValues = @("11","0.00002","0.00005","0.00003","0.00001","0.00003","0.00003","0.00004","0.00038","100.0")
$insertBefore = $t.Rows.Item(4)
foreach ($v in $newValues) {
    $t.Rows.Add($insertBefore) | Out-Null
}
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# After inserting 10 rows, the old row 34 (tab-separated values) is now row 44,
# the old row 35 is now row 45, and the old empty row 36 is now row 46.

# 5) Row 44 (was 34): collapse tab-separated run down to a single "100"
$t.Cell(44,1).Range.Text = "100"

# 6) Row 45 (was 35): collapse tab-separated run down to a single "0"
$t.Cell(45,1).Range.Text = "0"

# 7) Row 46 (was 36, empty): set text to "63"
$t.Cell(46,1).Range.Text = "63"
